# Add data for 2024-07-08: update 2024 (column K) year-to-date crime
# counts (and a couple of small 2016 corrections in column C) across the
# Citywide Totals, By Neighborhood, and per-neighborhood detail sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 5713
$ws.Range("K2").Value = 4051
$ws.Range("K3").Value = 4158
$ws.Range("K4").Value = 835
$ws.Range("K5").Value = 295
$ws.Range("K6").Value = 4652
$ws.Range("C7").Value = 28394
$ws.Range("K7").Value = 13991

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 120
$ws.Range("K5").Value = 33
$ws.Range("K7").Value = 406
$ws.Range("K8").Value = 954
$ws.Range("K9").Value = 57
$ws.Range("K11").Value = 275
$ws.Range("K16").Value = 42
$ws.Range("K18").Value = 95
$ws.Range("K19").Value = 428
$ws.Range("K20").Value = 311
$ws.Range("K23").Value = 140
$ws.Range("K25").Value = 64
$ws.Range("K27").Value = 138
$ws.Range("K29").Value = 740
$ws.Range("K31").Value = 153
$ws.Range("K33").Value = 581
$ws.Range("K36").Value = 176
$ws.Range("K37").Value = 476
$ws.Range("K42").Value = 497
$ws.Range("K43").Value = 124
$ws.Range("K44").Value = 128
$ws.Range("K47").Value = 80
$ws.Range("K48").Value = 180
$ws.Range("K49").Value = 80
$ws.Range("K51").Value = 170
$ws.Range("K52").Value = 381
$ws.Range("K53").Value = 187
$ws.Range("K55").Value = 157
$ws.Range("C63").Value = 278
$ws.Range("K63").Value = 43
$ws.Range("K65").Value = 322
$ws.Range("K68").Value = 35
$ws.Range("K71").Value = 43
$ws.Range("K73").Value = 128
$ws.Range("K74").Value = 15
$ws.Range("K75").Value = 46
$ws.Range("K76").Value = 199
$ws.Range("K77").Value = 98
$ws.Range("K78").Value = 166
$ws.Range("K79").Value = 364
$ws.Range("K83").Value = 299
$ws.Range("K84").Value = 101
$ws.Range("K85").Value = 630
$ws.Range("K88").Value = 161
$ws.Range("K89").Value = 197
$ws.Range("K90").Value = 129
$ws.Range("K91").Value = 153
$ws.Range("K92").Value = 49
$ws.Range("K93").Value = 49
$ws.Range("K94").Value = 176
$ws.Range("K95").Value = 237
$ws.Range("K99").Value = 241
$ws.Range("C101").Value = 28394
$ws.Range("K101").Value = 13991

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 145
$ws.Range("K3").Value = 133
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 406

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 210
$ws.Range("K7").Value = 630

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 100
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 381

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 268
$ws.Range("K3").Value = 284
$ws.Range("K5").Value = 27
$ws.Range("K6").Value = 321
$ws.Range("K7").Value = 954

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 299

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 219
$ws.Range("K7").Value = 581

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 84
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 237

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 159
$ws.Range("K6").Value = 142
$ws.Range("K7").Value = 476

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 95
$ws.Range("K7").Value = 322

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 98
$ws.Range("K7").Value = 241

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 54
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 210
$ws.Range("K3").Value = 264
$ws.Range("K5").Value = 22
$ws.Range("K7").Value = 740

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 42
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 130
$ws.Range("K6").Value = 129
$ws.Range("K7").Value = 428

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 162
$ws.Range("K6").Value = 178
$ws.Range("K7").Value = 497

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 51
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 123
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 364

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 108
$ws.Range("K3").Value = 95
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 311

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 32
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 35
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 25
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 15
